# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "NIT-9013624874" estado de cuenta sheet listed three worker/period
# rows (16-18). The first one (NOLBERTO MATEUS RODRIGUEZ, period 2502) is
# removed completely, and the duplicate extra period row for FIDEL ANDRES
# CARVAJAL HERNANDEZ (period 2503) is removed as well, leaving only his
# 2504 period row. The summary cells (total mora, worker/period counts)
# are updated to match the now-single remaining row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "NOLBERTO MATEUS RODRIGUEZ" detail row (was row 16). This
# shifts the FIDEL ANDRES CARVAJAL HERNANDEZ / period 2504 row up from 17
# to 16, and the period-2503 row up from 18 to 17.
$ws.Rows("16:16").Delete()

# Remove the now-duplicate FIDEL ANDRES CARVAJAL HERNANDEZ period-2503 row
# (now sitting at row 17), leaving only the period-2504 row at row 16.
$ws.Rows("17:17").Delete()

# Update the "VALOR MORA" total to reflect the single remaining row.
$ws.Range("E11").Value = 16000

# Update "Cant. Trabajadores" (now a single worker) and "Cant. Periodos"
# (now a single period) counters.
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1
